# "Generate Report for Handoff"
# The 85f214c1-4de6-4f4d-b17e-8e6646aeb0bb.md file has been handed off for
# localization (zh-cn / de-de), moving it from "In Translation" to
# "Ready for handoff", with the Priority downgraded from "ht" to "mt" and
# fresh handoff timestamps recorded.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the 85f214c1-... file ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 18:14:04"

# --- zh-cn sheet: row 3 is the 85f214c1-... file (Source File Name) ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-18 18:13:57"

# --- de-de sheet: row 3 is the 85f214c1-... file (Source File Name) ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-18 18:14:04"

# Re-fit the Status-type columns now that "Ready for handoff" is longer
# than the previous "In Translation" text that sized them.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
